$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the picture placeholder shape (the one created from the "pic" layout
# placeholder, idx=15) that needs to move out of its placeholder and into an
# explicit, fixed location.
$ph = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Name -eq "Picture Placeholder 3") {
        $ph = $candidate
    }
}

# The COM host here hands out the lowest currently-unused shape Id to any
# newly created shape. Real PowerPoint instead keeps a monotonically
# increasing counter across the whole editing session, so the picture that
# eventually gets re-created lands on Id 13 (one past the highest Id ever
# used on this slide, which is 12 - IBMLogo). Temporarily occupy every free
# Id slot below 13 with throwaway shapes so that when we duplicate the
# picture it is forced onto Id 13, matching the real authoring trace; then
# remove the throwaway shapes again.
$usedIds = @()
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $usedIds += $s.Shapes.Item($i).Id
}
$fillers = @()
for ($candidateId = 1; $candidateId -lt 13; $candidateId++) {
    if ($usedIds -notcontains $candidateId) {
        $fillers += $s.Shapes.AddShape(1, 0, 0, 1, 1)
    }
}

# Duplicate the placeholder picture so the copy keeps a valid embedded-image
# reference, then drop the original placeholder shape entirely (mirrors the
# diff removing the old <p:pic> that used the "pic" placeholder).
$dupRange = $ph.Duplicate()
$newPic = $dupRange.Item(1)
$ph.Delete()

# Release the Id-filler shapes now that the new picture has claimed Id 13.
foreach ($filler in $fillers) {
    $filler.Delete()
}

# Rename and reposition the freestanding picture at its intended, explicit
# location/size (matches the new <a:xfrm> in the diff).
$newPic.Name = "Picture 12"
$newPic.Left = 7502051 / 12700.0
$newPic.Top = 570271 / 12700.0
$newPic.Width = 3952530 / 12700.0
$newPic.Height = 2438400 / 12700.0

# Drop the placeholder crop (<a:srcRect t="8872" b="8872"/> -> none) and force
# explicit rect geometry to be written out for the now free-standing picture.
$newPic.PictureFormat.CropTop = 0
$newPic.PictureFormat.CropBottom = 0
$newPic.PictureFormat.CropLeft = 0
$newPic.PictureFormat.CropRight = 0
$newPic.AutoShapeType = 1
